$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# ---------------------------------------------------------------
# Sheet view: scroll position and selection
# ---------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 2
$ws.Range("I9").Select()

# ---------------------------------------------------------------
# Column N width
# ---------------------------------------------------------------
$ws.Columns.Item(14).ColumnWidth = 60.42578125

# ---------------------------------------------------------------
# Row 3 height + wrap text for the 3rd header cell group (G3:I3)
# ---------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 54
$ws.Range("G3:I3").WrapText = $true

# ---------------------------------------------------------------
# Row 8: new reporting period data
# ---------------------------------------------------------------
$ws.Range("A8").Value = 2023
$ws.Range("B8").Value = 44927
$ws.Range("C8").Value = 45016
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = ""
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = "Departamento de Recursos Humanos (UPP)"
$ws.Range("L8").Value = 45026
$ws.Range("M8").Value = 45026
$ws.Range("N8").Value = "Para este periodo no se cuenta con personal jubilado o pensionado"

# Give the whole row a thin box border and left alignment ...
$rowRange = $ws.Range("A8:N8")
$rowRange.Borders.LineStyle = 1
$rowRange.HorizontalAlignment = -4131

# ... then restore/re-apply the date number format on the date cells
# (must happen *after* the border/alignment pass so the engine keeps
# reusing the builtin date format instead of minting a new one)
$ws.Range("B8").NumberFormat = "mm-dd-yy"
$ws.Range("C8").NumberFormat = "mm-dd-yy"
$ws.Range("L8").NumberFormat = "mm-dd-yy"
$ws.Range("M8").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------
# Data validation ranges shrink from row 201 to row 153
# ---------------------------------------------------------------
$ws.Range("D8:D201").Validation.Delete()
$ws.Range("J8:J201").Validation.Delete()

$v1 = $ws.Range("D8:D153").Validation
$v1.Add(3, 1, 1, "Hidden_13")
$v1.IgnoreBlank = $true
$v1.InCellDropdown = $true
$v1.ShowInput = $false
$v1.ShowError = $true

$v2 = $ws.Range("J8:J153").Validation
$v2.Add(3, 1, 1, "Hidden_29")
$v2.IgnoreBlank = $true
$v2.InCellDropdown = $true
$v2.ShowInput = $false
$v2.ShowError = $true
